$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FAOSTAT IPCC")

# The table had three trailing duplicate rows (21/22 dup 19/20, etc.) -
# drop them first so the classification table ends at row 19.
$ws.Rows("20:22").Delete()

# Re-map the IPCC codes in column A to the more granular UN
# classification (the source data now distinguishes sub-items such as
# 3.A.2.a/b/c instead of the single 3.A.2 bucket, etc.).
$ws.Range("A3").Value = "3.C.7"
$ws.Range("A2").Value = "9.A"
$ws.Range("A7").Value = "3.A.2.a"
$ws.Range("A8").Value = "3.A.2.b"
$ws.Range("A9").Value = "3.A.2.c"
$ws.Range("A16").Value = "3.B.1.a"
$ws.Range("A17").Value = "3.B.1.b"
$ws.Range("A18").Value = "3.B.6.a"
$ws.Range("A19").Value = "3.B.6.b"
$ws.Range("A12").Value = "3.C.1.a"
$ws.Range("A13").Value = "3.C.1.b"
$ws.Range("A14").Value = "3.C.1.c"
$ws.Range("A4").Value = "3.C.8"
$ws.Range("A5").Value = "3.A.1"
$ws.Range("A6").Value = "3.C.6"
$ws.Range("A10").Value = "3.C.4"
$ws.Range("A11").Value = "3.C.5"
$ws.Range("A15").Value = "3.B.2"

# Realign the "IPCC" driver label (col B) for each row with its new code.
$ws.Range("B2").Value = "Animal Stocks"
$ws.Range("B3").Value = "Agriculture - Rice cultivation (Emissions CH4)"
$ws.Range("B4").Value = "Agriculture - Crop residues (Emissions N2O)"
$ws.Range("B5").Value = "Agriculture - Enteric fermentation (Emissions CH4)"
$ws.Range("B6").Value = "Agriculture - Manure left on pasture (Emissions N2O)"
$ws.Range("B7").Value = "Agriculture - Emissions (N2O) (Manure applied)"
$ws.Range("B8").Value = "Agriculture - Manure management (Emissions N2O)"
$ws.Range("B9").Value = "Agriculture - Manure management (Emissions CH4)"
$ws.Range("B10").Value = "Agriculture - Synthetic Fertilizers"
$ws.Range("B11").Value = "Agriculture - Drained organic soils (N2O)"
$ws.Range("B12").Value = "Agriculture - Burning crop residues (Emissions CH4)"
$ws.Range("B13").Value = "Agriculture - Burning crop residues (Emissions N2O)"
$ws.Range("B14").Value = "Agriculture - Savanna fires"
$ws.Range("B15").Value = "LULUCF - Drained organic soils (CO2)"
$ws.Range("B16").Value = "LULUCF - Forestland"
$ws.Range("B17").Value = "LULUCF - Net Forest conversion"
$ws.Range("B18").Value = "LULUCF - Fires in organic soils"
$ws.Range("B19").Value = "LULUCF - Forest fires"

# ... and the short "Driver" label (col C) to match.
$ws.Range("C2").Value = "Stocks"
$ws.Range("C3").Value = "Rice cultivation (Emissions CH4)"
$ws.Range("C4").Value = "Crop residues (Emissions N2O)"
$ws.Range("C5").Value = "Enteric fermentation (Emissions CH4)"
$ws.Range("C6").Value = "Manure left on pasture (Emissions N2O)"
$ws.Range("C7").Value = "Emissions (N2O) (Manure applied)"
$ws.Range("C8").Value = "Manure management (Emissions N2O)"
$ws.Range("C9").Value = "Manure management (Emissions CH4)"
$ws.Range("C10").Value = "Synthetic Fertilizers"
$ws.Range("C11").Value = "Drained organic soils (N2O)"
$ws.Range("C12").Value = "Burning crop residues (Emissions CH4)"
$ws.Range("C13").Value = "Burning crop residues (Emissions N2O)"
$ws.Range("C14").Value = "Savanna fires"
$ws.Range("C15").Value = "Drained organic soils (CO2)"
$ws.Range("C16").Value = "Forestland"
$ws.Range("C17").Value = "Net Forest conversion"
$ws.Range("C18").Value = "Fires in organic soils"
$ws.Range("C19").Value = "Forest fires"

# Column B now holds the longer "Agriculture - ..." / "LULUCF - ..."
# driver labels, so widen it to fit them.
$ws.Columns("B:B").ColumnWidth = 42.86

# Restore the recorded view state (scroll position + active selection).
$ws.Range("D25").Select()
